$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A15").Value = "Pulse Measurement"
$ws.Range("D15").Value = "PF0 / PWM0"

$ws.Range("A16").Select()
